# Applies the regression/constants-evaluation refresh described in the commit:
# "calorimetry : gui regression tests : NMR, spectrophotometry"
# Updates cached numeric results across several computed sheets in the workbook.

$wb = $excel.ActiveWorkbook

# --- constants_evaluated: B4/C4/B5/C5 hold numeric-looking text (shared strings) ---
# Force text storage (NumberFormat "@") while assigning, then restore the Normal
# style so no stray number-format style is left behind on the cells.
$wsConstEval = $wb.Worksheets.Item("constants_evaluated")
$rngConstEval = $wsConstEval.Range("B4:C5")
$rngConstEval.NumberFormat = "@"
$wsConstEval.Range("B4").Value = "4.571025390625"
$wsConstEval.Range("C4").Value = "0.108717276173831"
$wsConstEval.Range("B5").Value = "6.46953125"
$wsConstEval.Range("C5").Value = "0.260200371728513"
$rngConstEval.Style = "Normal"

# --- correlation_matrix: refreshed regression output values ---
$wsCorr = $wb.Worksheets.Item("correlation_matrix")
$wsCorr.Range("B2").Value = 0.848847310230273
$wsCorr.Range("A3").Value = 0.848847310230273

# --- adj_r_squared: refreshed regression output values ---
$wsAdjR2 = $wb.Worksheets.Item("adj_r_squared")
$wsAdjR2.Range("A2").Value = 0.999944167667562

# --- mol_ext_coefficients_calc: refreshed regression output values ---
$wsMolExt = $wb.Worksheets.Item("mol_ext_coefficients_calc")
$wsMolExt.Range("C2").Value = 3.62851768723924
$wsMolExt.Range("D2").Value = 2964.36394125175
$wsMolExt.Range("E2").Value = 4768.81513735542
$wsMolExt.Range("F2").Value = 5055.18778225329
$wsMolExt.Range("G2").Value = 423108.972180095
$wsMolExt.Range("C3").Value = 22.5538789403126
$wsMolExt.Range("D3").Value = 8048.51242828701
$wsMolExt.Range("E3").Value = 6303.21573312717
$wsMolExt.Range("F3").Value = 3311.39032240926
$wsMolExt.Range("G3").Value = -639364.45312061
$wsMolExt.Range("C4").Value = 2.2557199599467
$wsMolExt.Range("D4").Value = 70.1880926029722
$wsMolExt.Range("E4").Value = 12.6654507470743
$wsMolExt.Range("F4").Value = 164.100837409176
$wsMolExt.Range("G4").Value = 109207.061021345
$wsMolExt.Range("C5").Value = 1.80524215495195
$wsMolExt.Range("D5").Value = 56.1712028941527
$wsMolExt.Range("E5").Value = 10.1361011145313
$wsMolExt.Range("F5").Value = 131.329134207315
$wsMolExt.Range("G5").Value = 87397.9011910678

# --- equilibrium_concentrations: refreshed regression output values ---
$wsEquilConc = $wb.Worksheets.Item("equilibrium_concentrations")
$wsEquilConc.Range("A2").Value = 0.0000000657079255915789
$wsEquilConc.Range("B2").Value = 0.000353335365591123
$wsEquilConc.Range("C2").Value = 0.000000864629911626739
$wsEquilConc.Range("D2").Value = 0.00000000000449732084461406
$wsEquilConc.Range("E2").Value = 0.000000230346831802918
$wsEquilConc.Range("A3").Value = 0.0000839967070662358
$wsEquilConc.Range("B3").Value = 0.0000857807350894835
$wsEquilConc.Range("C3").Value = 0.000268335056694733
$wsEquilConc.Range("D3").Value = 0.00000178420821621364
$wsEquilConc.Range("E3").Value = 0.000000000180192926758746
$wsEquilConc.Range("A4").Value = 0.000347583618662726
$wsEquilConc.Range("B4").Value = 0.000024670213965408
$wsEquilConc.Range("C4").Value = 0.000319343151640457
$wsEquilConc.Range("D4").Value = 0.00000878663865680689
$wsEquilConc.Range("E4").Value = 0.0000000000435452411209539
$wsEquilConc.Range("A5").Value = 0.000665672088035321
$wsEquilConc.Range("B5").Value = 0.0000131491284715596
$wsEquilConc.Range("C5").Value = 0.00032597380835523
$wsEquilConc.Range("D5").Value = 0.0000171770631735829
$wsEquilConc.Range("E5").Value = 0.0000000000227373398350435
$wsEquilConc.Range("A6").Value = 0.00133310373218991
$wsEquilConc.Range("B6").Value = 0.00000636835440170605
$wsEquilConc.Range("C6").Value = 0.000316167015005058
$wsEquilConc.Range("D6").Value = 0.0000333646338094201
$wsEquilConc.Range("E6").Value = 0.0000000000113536644740305
$wsEquilConc.Range("A7").Value = 0.00264712583306349
$wsEquilConc.Range("B7").Value = 0.00000294410784641187
$wsEquilConc.Range("C7").Value = 0.000290237611670279
$wsEquilConc.Range("D7").Value = 0.0000608182805040659
$wsEquilConc.Range("E7").Value = 0.00000000000571775330636466
$wsEquilConc.Range("A8").Value = 0.00664397774570443
$wsEquilConc.Range("B8").Value = 0.000000935641668413102
$wsEquilConc.Range("C8").Value = 0.000231506460134426
$wsEquilConc.Range("D8").Value = 0.000121757898265341
$wsEquilConc.Range("E8").Value = 0.00000000000227809500026514
$wsEquilConc.Range("A9").Value = 0.0134044254224832
$wsEquilConc.Range("B9").Value = 0.00000034362691489237
$wsEquilConc.Range("C9").Value = 0.00017153816754401
$wsEquilConc.Range("D9").Value = 0.000182018205581655
$wsEquilConc.Range("E9").Value = 0.00000000000112915041169726

# --- absorbance_calc_abs_errors: refreshed regression output values ---
$wsAbsErr = $wb.Worksheets.Item("absorbance_calc_abs_errors")
$wsAbsErr.Range("C2").Value = 1.14899994954435
$wsAbsErr.Range("D2").Value = 1.5433261305592
$wsAbsErr.Range("E2").Value = 1.64171769433111
$wsAbsErr.Range("F2").Value = 1.68274593717371
$wsAbsErr.Range("G2").Value = 1.70012665074334
$wsAbsErr.Range("H2").Value = 1.70987231378232
$wsAbsErr.Range("I2").Value = 1.74640288825853
$wsAbsErr.Range("J2").Value = 1.78782732674432
$wsAbsErr.Range("C3").Value = 2.70199995089303
$wsAbsErr.Range("D3").Value = 2.38946851601394
$wsAbsErr.Range("E3").Value = 2.24835480918796
$wsAbsErr.Range("F3").Value = 2.23239307233761
$wsAbsErr.Range("G3").Value = 2.18466740939406
$wsAbsErr.Range("H3").Value = 2.11421833415463
$wsAbsErr.Range("I3").Value = 2.01979962465506
$wsAbsErr.Range("J3").Value = 1.9890621527603
$wsAbsErr.Range("C4").Value = -0.0000000504556501024211
$wsAbsErr.Range("D4").Value = 0.000326130559196391
$wsAbsErr.Range("E4").Value = -0.000282305668892313
$wsAbsErr.Range("F4").Value = -0.00225406282629392
$wsAbsErr.Range("G4").Value = -0.000873349256657052
$wsAbsErr.Range("H4").Value = 0.00587231378232156
$wsAbsErr.Range("I4").Value = -0.00359711174146571
$wsAbsErr.Range("J4").Value = 0.000827326744315293
$wsAbsErr.Range("C5").Value = -0.0000000491069669372735
$wsAbsErr.Range("D5").Value = 0.000468516013935272
$wsAbsErr.Range("E5").Value = -0.00364519081203563
$wsAbsErr.Range("F5").Value = 0.00439307233761443
$wsAbsErr.Range("G5").Value = -0.00133259060594426
$wsAbsErr.Range("H5").Value = 0.00021833415462913
$wsAbsErr.Range("I5").Value = -0.000200375344935111
$wsAbsErr.Range("J5").Value = 0.0000621527603048566

# --- absorbance_calc_rel_errors: refreshed regression output values ---
$wsRelErr = $wb.Worksheets.Item("absorbance_calc_rel_errors")
$wsRelErr.Range("C2").Value = 1.14899994954435
$wsRelErr.Range("D2").Value = 1.5433261305592
$wsRelErr.Range("E2").Value = 1.64171769433111
$wsRelErr.Range("F2").Value = 1.68274593717371
$wsRelErr.Range("G2").Value = 1.70012665074334
$wsRelErr.Range("H2").Value = 1.70987231378232
$wsRelErr.Range("I2").Value = 1.74640288825853
$wsRelErr.Range("J2").Value = 1.78782732674432
$wsRelErr.Range("C3").Value = 2.70199995089303
$wsRelErr.Range("D3").Value = 2.38946851601394
$wsRelErr.Range("E3").Value = 2.24835480918796
$wsRelErr.Range("F3").Value = 2.23239307233761
$wsRelErr.Range("G3").Value = 2.18466740939406
$wsRelErr.Range("H3").Value = 2.11421833415463
$wsRelErr.Range("I3").Value = 2.01979962465506
$wsRelErr.Range("J3").Value = 1.9890621527603
$wsRelErr.Range("C4").Value = -0.000000043912663274518
$wsRelErr.Range("D4").Value = 0.000211361347502522
$wsRelErr.Range("E4").Value = -0.000171927934769984
$wsRelErr.Range("F4").Value = -0.00133772274557503
$wsRelErr.Range("G4").Value = -0.000513432837540889
$wsRelErr.Range("H4").Value = 0.00344619353422627
$wsRelErr.Range("I4").Value = -0.00205549242369469
$wsRelErr.Range("J4").Value = 0.000462969638676717
$wsRelErr.Range("C5").Value = -0.0000000181743030855934
$wsRelErr.Range("D5").Value = 0.000196113861002625
$wsRelErr.Range("E5").Value = -0.00161864600889681
$wsRelErr.Range("F5").Value = 0.00197175598636195
$wsRelErr.Range("G5").Value = -0.00060960229000195
$wsRelErr.Range("H5").Value = 0.000103280110988235
$wsRelErr.Range("I5").Value = -0.0000991957153144113
$wsRelErr.Range("J5").Value = 0.0000312482455026931

Write-Output "edits applied"
